$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 13,12
$data[0,0] = 'nurse_adult_staff_codes'
$data[0,1] = 'nurse_children_staff_codes'
$data[0,2] = 'nurse_maternity_staff_codes'
$data[0,3] = 'nurse_community_mental_health_staff_codes'
$data[0,4] = 'nurse_other_mental_health_staff_codes'
$data[0,5] = 'nurse_community_learning_difficulties_staff_codes'
$data[0,6] = 'nurse_other_learning_difficulties_staff_codes'
$data[0,7] = 'nurse_community_service_staff_codes'
$data[0,8] = 'nurse_education_staff_codes'
$data[0,9] = 'nurse_school_nursing_staff_codes'
$data[0,10] = 'nurse_neonatal_staff_codes'
$data[0,11] = 'nurse_learner_staff_codes'
$data[1,0] = 'NAA'
$data[1,1] = 'NAB'
$data[1,2] = 'NAC'
$data[1,3] = 'NAD'
$data[1,4] = 'NAE'
$data[1,5] = 'NAF'
$data[1,6] = 'NAG'
$data[1,7] = 'NAH'
$data[1,8] = 'NAJ'
$data[1,9] = 'NAK'
$data[1,10] = 'NAL'
$data[1,11] = 'P2A'
$data[2,0] = 'NCA'
$data[2,1] = 'NCB'
$data[2,2] = 'NCC'
$data[2,3] = 'NCD'
$data[2,4] = 'NCE'
$data[2,5] = 'NCF'
$data[2,6] = 'NCG'
$data[2,7] = 'NCH'
$data[2,8] = 'NCJ'
$data[2,9] = 'NCK'
$data[2,10] = 'NCL'
$data[2,11] = 'P3A'
$data[3,0] = 'N0A'
$data[3,1] = 'N0B'
$data[3,2] = 'N0C'
$data[3,3] = 'N0D'
$data[3,4] = 'N0E'
$data[3,5] = 'N0E'
$data[3,6] = 'N0G'
$data[3,7] = 'NEH'
$data[3,8] = 'N0J'
$data[3,9] = 'N0K'
$data[3,10] = 'N0L'
$data[3,11] = 'P1A'
$data[4,0] = 'N6A'
$data[4,1] = 'N1B'
$data[4,2] = 'N1C'
$data[4,3] = 'N4D'
$data[4,4] = 'N6E'
$data[4,5] = 'N0F'
$data[4,6] = 'N6G'
$data[4,7] = 'N0H'
$data[4,8] = 'N1J'
$data[4,9] = 'NBK'
$data[4,10] = 'N1L'
$data[4,11] = 'P2B'
$data[5,0] = 'N7A'
$data[5,1] = 'N6B'
$data[5,2] = 'N2C'
$data[5,3] = 'N5D'
$data[5,4] = 'N7E'
$data[5,5] = 'N4F'
$data[5,6] = 'N7G'
$data[5,7] = 'N1H'
$data[5,8] = 'N3J'
$data[5,9] = 'N6K'
$data[5,10] = 'N2L'
$data[5,11] = 'P2C'
$data[6,1] = 'N7B'
$data[6,2] = 'N6C'
$data[6,3] = 'N6D'
$data[6,5] = 'N5F'
$data[6,7] = 'N3H'
$data[6,8] = 'N6J'
$data[6,9] = 'N7K'
$data[6,10] = 'N6L'
$data[6,11] = 'P3C'
$data[7,2] = 'N7C'
$data[7,3] = 'N7D'
$data[7,5] = 'N6F'
$data[7,7] = 'N4H'
$data[7,8] = 'N7J'
$data[7,10] = 'N7L'
$data[7,11] = 'P2D'
$data[8,5] = 'N7F'
$data[8,7] = 'N5H'
$data[8,11] = 'P3D'
$data[9,7] = 'N6H'
$data[9,11] = 'P1D'
$data[10,7] = 'N7H'
$data[10,11] = 'P2E'
$data[11,11] = 'P3E'
$data[12,11] = 'P1E'

$ws.Range("AA1:AL13").Value = $data

$ws.Range("AF6").Select()
